$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E are treated as Text so numeric-looking values
# (e.g. "63.234.97", "0.109") are preserved exactly as strings, matching
# the original inline-string cell content, instead of being coerced to numbers.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.234.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.11%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.238.27'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.84%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.18'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.39%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.10'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.62%  '

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.233.51'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.91%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.10%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.147'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.41%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.35'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.69%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.55%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.23%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.18'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.05%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.775.08'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.96%  '

# Row 16
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.28%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.238.36'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.04%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.284.73'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.13%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.77'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.95%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '473.18'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.08%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.11'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.84%  '

# Row 22
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.36%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.93'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.03%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.42'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.29%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.16'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.89%  '

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.06%  '

# Row 27
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.73'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.23%  '

# Row 28
$ws.Range('B28').Value = 'NEARProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.39'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.07%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.07'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.50%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.38%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.39'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.65%  '

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.04%  '

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.33%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.52'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.69%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.08'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.81%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.90'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.45%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.64'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.00%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0705'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.86%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0392'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.73%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '421.35'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.46%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.35'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.21%  '

# Row 42
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.958.94'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.20%  '

# Row 43
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.73'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.82%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.109'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -9.31%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.266'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.23%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.53%  '

# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.06%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.34'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.65%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.85'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.24%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.76%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.21'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.70%  '
